$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-05-26 Monday" "2025-05-27 Tuesday"

Replace-Text "48÷2=" "82÷3="
Replace-Text "23÷8=" "39÷9="
Replace-Text "37÷9=" "36÷2="
Replace-Text "18÷4=" "64÷6="
Replace-Text "43÷6=" "22÷6="
Replace-Text "36÷7=" "51÷6="
Replace-Text "37÷3=" "54÷7="
Replace-Text "35÷3=" "93÷9="
Replace-Text "40÷6=" "77÷5="
Replace-Text "79÷4=" "56÷7="
Replace-Text "69÷8=" "46÷9="
Replace-Text "79÷6=" "36÷8="
Replace-Text "59÷5=" "72÷8="
Replace-Text "99÷2=" "41÷2="
Replace-Text "82÷6=" "87÷6="
Replace-Text "60÷3=" "86÷3="
Replace-Text "81÷6=" "57÷4="
Replace-Text "42÷6=" "49÷5="
Replace-Text "62÷6=" "30÷2="
Replace-Text "17÷3=" "71÷9="
Replace-Text "75÷6=" "43÷5="
Replace-Text "84÷2=" "90÷5="
Replace-Text "89÷4=" "96÷8="
Replace-Text "20÷4=" "53÷7="
Replace-Text "11÷5=" "23÷2="
